# 自动更新Excel文件脚本
# For every data row, the "剩余" (remaining days) counter in column E counts
# down by 1 each run. Once it reaches 1, the cycle restarts: E is reset to
# the row's total-days value (column D) and the start date in column F is
# bumped to the new cycle's start date (today, as YYYYMMDD). Rows whose F
# value isn't a well-formed 8-digit YYYYMMDD date (e.g. a bad/typo'd date)
# are left untouched, matching the source system's behaviour.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newStartDate = 20260223
$firstRow = 2
$lastRow = 99

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $dCell = $ws.Cells.Item($r, 4)
    $eCell = $ws.Cells.Item($r, 5)
    $fCell = $ws.Cells.Item($r, 6)

    $dVal = $dCell.Value2
    $eVal = $eCell.Value2
    $fVal = $fCell.Value2

    if ($dVal -eq $null -or $eVal -eq $null) {
        continue
    }

    $fText = [string]$fVal
    $isValidDate = ($fVal -ne $null) -and ($fText.Length -eq 8)

    if (-not $isValidDate) {
        continue
    }

    if ($eVal -le 1) {
        $eCell.Value2 = $dVal
        $fCell.Value2 = $newStartDate
    } else {
        $eCell.Value2 = $eVal - 1
    }
}
